# Insert two new weekly observations at the top of the "Ají" data block
# (row 124), pushing the existing 124:193 block down to 126:195 and growing
# the sheet from A1:R193 to A1:R195.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("124:125").Insert()

# New row 124
$ws.Cells.Item(124, 4).Value2  = 44529              # D - Fecha
$ws.Cells.Item(124, 8).Value2  = "Inferno"           # H - Variedad
$ws.Cells.Item(124, 9).Value2  = "Primera"           # I - Calidad
$ws.Cells.Item(124, 10).Value2 = 43                  # J - Volumen
$ws.Cells.Item(124, 11).Value2 = 22000               # K - Precio minimo
$ws.Cells.Item(124, 12).Value2 = 25000               # L - Precio maximo
$ws.Cells.Item(124, 13).Value2 = 23535               # M - Precio promedio
$ws.Cells.Item(124, 14).Value2 = "$/caja 12 kilos"   # N - Unidad
$ws.Cells.Item(124, 15).Value2 = "Región de Arica y Parinacota"  # O - Origen
$ws.Cells.Item(124, 16).Value2 = 1961                # P - Precio promedio $/kilo
$ws.Cells.Item(124, 17).Value2 = 12                  # Q - Kilos por unidad

# New row 125
$ws.Cells.Item(125, 4).Value2  = 44529               # D - Fecha
$ws.Cells.Item(125, 8).Value2  = "Inferno"           # H - Variedad
$ws.Cells.Item(125, 9).Value2  = "Primera"           # I - Calidad
$ws.Cells.Item(125, 10).Value2 = 25                  # J - Volumen
$ws.Cells.Item(125, 11).Value2 = 28000               # K - Precio minimo
$ws.Cells.Item(125, 12).Value2 = 30000               # L - Precio maximo
$ws.Cells.Item(125, 13).Value2 = 28960               # M - Precio promedio
$ws.Cells.Item(125, 14).Value2 = "$/caja 15 kilos"   # N - Unidad
$ws.Cells.Item(125, 15).Value2 = "Provincia de Huasco"  # O - Origen
$ws.Cells.Item(125, 16).Value2 = 1931                # P - Precio promedio $/kilo
$ws.Cells.Item(125, 17).Value2 = 15                  # Q - Kilos por unidad

# The other (constant) columns for the two new rows mirror the rest of the block
$ws.Cells.Item(124, 1).Value2 = 9
$ws.Cells.Item(124, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(124, 3).Value2 = "Metropolitana"
$ws.Cells.Item(124, 5).Value2 = 13
$ws.Cells.Item(124, 6).Value2 = 100112021
$ws.Cells.Item(124, 7).Value2 = "Ají"
$ws.Cells.Item(124, 18).Value2 = "Hortaliza"

$ws.Cells.Item(125, 1).Value2 = 9
$ws.Cells.Item(125, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(125, 3).Value2 = "Metropolitana"
$ws.Cells.Item(125, 5).Value2 = 13
$ws.Cells.Item(125, 6).Value2 = 100112021
$ws.Cells.Item(125, 7).Value2 = "Ají"
$ws.Cells.Item(125, 18).Value2 = "Hortaliza"
